$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 62, shifting rows 62:90 down to 63:91
$ws.Rows("62:62").Insert()

# Populate the new row 62 with data (same as old row 62 but with updated D/N/O/P/S)
$ws.Cells.Item(62, 1).Value = 1
$ws.Cells.Item(62, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(62, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(62, 4).Value = 44636
$ws.Cells.Item(62, 5).Value = 15
$ws.Cells.Item(62, 6).Value = "Fruta"
$ws.Cells.Item(62, 7).Value = 100102
$ws.Cells.Item(62, 8).Value = "Cítricos"
$ws.Cells.Item(62, 9).Value = 100102004
$ws.Cells.Item(62, 10).Value = "Mandarina"
$ws.Cells.Item(62, 11).Value = "Murcott"
$ws.Cells.Item(62, 12).Value = "Segunda"
$ws.Cells.Item(62, 13).Value = 270
$ws.Cells.Item(62, 14).Value = 19000
$ws.Cells.Item(62, 15).Value = 20000
$ws.Cells.Item(62, 16).Value = 19500
$ws.Cells.Item(62, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(62, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(62, 19).Value = 975
$ws.Cells.Item(62, 20).Value = 20
